$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$clusterNames = @{
    20 = "ECs"
    21 = "FAPs"
    22 = "M1"
    23 = "M2"
    24 = "Neutro"
    25 = "sCs"
}

# Each entry: row, sendingClusterIdx, targetClusterIdx, E..T values
$data = @(
    @(2,  20, 21, 3,1, 86.484492,          259.453476,         0.1133842074223504, 0.1133842074223504, 1, 0.3333333333333333, 0.2064483333333333, 0.619345,   0.3625320258231888, 0.3625320258231888, 17.85457923258,     160.69121309322,    0.04110540641318133, 0.04110540641318134),
    @(3,  20, 25, 3,1, 86.484492,          259.453476,         0.1133842074223504, 0.1133842074223504, 3, 1,                  0.363014,            1.089042,   0.6374679741768112, 0.6374679741768112, 31.395081378888,    282.5557324099921,  0.07227880100916907, 0.0722788010091691),
    @(4,  21, 21, 3,1, 157.8540903333333,  473.562271,         0.206952258224759,  0.206952258224759,  1, 0.3333333333333333, 0.2064483333333333, 0.619345,   0.3625320258231888, 0.3625320258231888, 32.58871385916611,  293.2984247324951,  0.07502682142290558, 0.07502682142290558),
    @(5,  21, 25, 3,1, 157.8540903333333,  473.562271,         0.206952258224759,  0.206952258224759,  3, 1,                  0.363014,            1.089042,   0.6374679741768112, 0.6374679741768112, 57.30324474826467,  515.7292027343821,  0.1319254368018535,  0.1319254368018535),
    @(6,  22, 21, 3,1, 107.1200406666667,  321.360122,         0.1404381367013169, 0.1404381367013169, 1, 0.3333333333333333, 0.2064483333333333, 0.619345,   0.3625320258231888, 0.3625320258231888, 22.11475386223222,  199.03278476009,    0.05091332220116233, 0.05091332220116233),
    @(7,  22, 25, 3,1, 107.1200406666667,  321.360122,         0.1404381367013169, 0.1404381367013169, 3, 1,                  0.363014,            1.089042,   0.6374679741768112, 0.6374679741768112, 38.88607444256933,  349.974669983124,   0.08952481450015454, 0.08952481450015455),
    @(8,  23, 21, 3,1, 104.2030356666667,  312.609107,         0.1366138406648433, 0.1366138406648433, 1, 0.3333333333333333, 0.2064483333333333, 0.619345,   0.3625320258231888, 0.3625320258231888, 21.51254304165722,  193.612887374915,   0.04952689241171196, 0.04952689241171196),
    @(9,  23, 25, 3,1, 104.2030356666667,  312.609107,         0.1366138406648433, 0.1366138406648433, 3, 1,                  0.363014,            1.089042,   0.6374679741768112, 0.6374679741768112, 37.82716078949933,  340.444447105494,   0.0870869482531313,  0.08708694825313132),
    @(10, 24, 21, 3,1, 36.59512433333333,  109.785373,         0.04797749367663938,0.04797749367663939,1, 0.3333333333333333, 0.2064483333333333, 0.619345,   0.3625320258231888, 0.3625320258231888, 7.555002426742778,  67.995021840685,    0.01739337797651131, 0.01739337797651131),
    @(11, 24, 25, 3,1, 36.59512433333333,  109.785373,         0.04797749367663938,0.04797749367663939,3, 1,                  0.363014,            1.089042,   0.6374679741768112, 0.6374679741768112, 13.28454246474067,  119.560882182666,   0.03058411570012808, 0.03058411570012808),
    @(12, 25, 21, 3,1, 270.499283,         811.497849,         0.3546340633100911, 0.3546340633100911, 1, 0.3333333333333333, 0.2064483333333333, 0.619345,   0.3625320258231888, 0.3625320258231888, 55.84412614321167,  502.597135288905,   0.1285662053977163,  0.1285662053977163),
    @(13, 25, 25, 3,1, 270.499283,         811.497849,         0.3546340633100911, 0.3546340633100911, 3, 1,                  0.363014,            1.089042,   0.6374679741768112, 0.6374679741768112, 98.195026718962,    883.755240470658,   0.2260678579123748,  0.2260678579123748)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $sendingIdx = $entry[1]
    $targetIdx = $entry[2]

    $ws.Cells.Item($r, 1).Value = $clusterNames[$sendingIdx]
    $ws.Cells.Item($r, 2).Value = "Gnas"
    $ws.Cells.Item($r, 3).Value = "Lhcgr"
    $ws.Cells.Item($r, 4).Value = $clusterNames[$targetIdx]

    for ($i = 3; $i -lt $entry.Length; $i++) {
        $col = $i + 2  # entry[3] -> column E (5)
        $ws.Cells.Item($r, $col).Value = $entry[$i]
    }
}
